$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("File Name", "Quantity"),
    @("motor_holder.STL", 2),
    @("bearing_block.STL", 4),
    @("wheel_mount.STL", 4),
    @("3_in_stand_offs.STL", 8),
    @("battery_holder.STL", 2),
    @("jetson_orin_nano_holder.STL", 1),
    @("9_inch_spacers_3prong.STL", 12),
    @("astra_holder.STL", 1),
    @("3.5_inch_cir_standoff.STL", 4),
    @("bottom_support_left.STL", 1),
    @("bottom_support_right.STL", 1),
    @("top_support_left.STL", 1),
    @("top_support_right.STL", 1),
    @("microphone_casing_back.STL", 2),
    @("microphone_casing.STL", 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A1:B1").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth = 23.17

$ws.Range("C4").Select() | Out-Null
